$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking") updates
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total") updates
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -12
$ws.Range("E12").Value = "68.0/140"
